$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 5 ("ZohoTitle") so a new "ProcessedFolderPath"
# / "ProcessedFolderPathEnd" pair can be inserted right after the existing
# "RequestFilePath" rows -- i.e. the generated file now gets moved to a
# Processed folder right after the request is created.
$ws.Range("A5:B6").EntireRow.Insert()

# Row 6 first (so its shared string "ProcessedFolderPath" is registered
# before "ProcessedFolderPathEnd"/the literal path text below it).
$ws.Cells.Item(6, 1).Value = "ProcessedFolderPath"

# Row 5: ProcessedFolderPathEnd / relative processed-folder path
$ws.Cells.Item(5, 1).Value = "ProcessedFolderPathEnd"
$ws.Cells.Item(5, 2).Value = "UiPath\HelpDeskTicketGeneration\Data\Processed\"

# Row 6: ProcessedFolderPath = CONCAT(DocumentsPath, ProcessedFolderPathEnd)
$ws.Cells.Item(6, 2).Formula = "=CONCAT(B2,B5)"

# Apply the column B data style (fill/quote-prefix format used by every
# other Value cell) down from row 4 onto the two new rows, now that their
# values are in place.
$ws.Range("B4").Copy()
$ws.Range("B5:B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A's longest entry is now "ProcessedFolderPathEnd" -- re-apply the
# best-fit auto width so the column keeps showing the key names in full.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(1).ColumnWidth = 20.75

$ws.Range("A6").Select()

$wb.Save()
